$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow edits, then restore protection
$ws.Unprotect()

# Update the confidential/disclosure text (date change 2021-06-09 -> 2021-06-10)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5411588762757222
$ws.Range("E2").Value = -0.004128929142248161

$ws.Range("D3").Value = 0.2511652708071809
$ws.Range("E3").Value = 0.008221549112938265

$ws.Range("D4").Value = 0.04986885455620273
$ws.Range("E4").Value = 0.002996628792608202

$ws.Range("D5").Value = 0.09886485273803498
$ws.Range("E5").Value = -0.007328691828508704

$ws.Range("D6").Value = 0.02892539606376871
$ws.Range("E6").Value = -0.01844748858447487

$ws.Range("D7").Value = 0.03001674955909055
$ws.Range("E7").Value = -0.01221995926680253

$ws.Range("E8").Value = -0.001644955009027504

# Restore sheet protection as it was originally
$ws.Protect()
